# Move the "x" mark in the "Essenciais" sheet from C4 ("Não realizado")
# to D4 ("Em realização") for the row "Em conjunto com o APF, preparar o
# seu Plano Pessoal de Formação." (row 4).
#
# All downstream totals/percentages on "Essenciais" (C17/D17, C19/D19)
# and on "Resultados" (C4/D4, C5/D5, C8/D8, C9/D9) are plain formulas and
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Essenciais")

$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "x"

# Leave the selection on C4, matching the saved view state.
$ws.Range("C4").Select()
